# Simulated Wild Card round and logged it
# Update the "R" (Road/Runningback? row 3) stats on both the OFF and DEF
# sheets to reflect the results of the simulated Wild Card round game.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 483
$wsOff.Range("C3").Value = 342
$wsOff.Range("D3").Value = 100
$wsOff.Range("E3").Value = 45
$wsOff.Range("G3").Value = 10

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 387
$wsDef.Range("C3").Value = 264
$wsDef.Range("D3").Value = 106
$wsDef.Range("E3").Value = 46
